# Update the generated multiplication-table answers to a new random set.
# Each cell is addressed by its (row, column) position in the single table
# so that duplicate values (e.g. "29x83=2407" appears twice in the source)
# are not conflated by a global find/replace.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "73×16=1168"; New = "59×44=2596" },
    @{ Row = 1;  Col = 2; Old = "17×61=1037"; New = "77×72=5544" },
    @{ Row = 1;  Col = 3; Old = "90×29=2610"; New = "55×93=5115" },
    @{ Row = 1;  Col = 4; Old = "24×83=1992"; New = "25×50=1250" },
    @{ Row = 1;  Col = 5; Old = "29×83=2407"; New = "68×32=2176" },

    @{ Row = 5;  Col = 1; Old = "37×41=1517"; New = "14×62=868"  },
    @{ Row = 5;  Col = 2; Old = "57×52=2964"; New = "93×56=5208" },
    @{ Row = 5;  Col = 3; Old = "76×68=5168"; New = "74×98=7252" },
    @{ Row = 5;  Col = 4; Old = "46×79=3634"; New = "12×83=996"  },
    @{ Row = 5;  Col = 5; Old = "52×54=2808"; New = "68×73=4964" },

    @{ Row = 10; Col = 1; Old = "87×63=5481"; New = "90×57=5130" },
    @{ Row = 10; Col = 2; Old = "83×60=4980"; New = "82×12=984"  },
    @{ Row = 10; Col = 3; Old = "34×87=2958"; New = "12×64=768"  },
    @{ Row = 10; Col = 4; Old = "93×41=3813"; New = "58×30=1740" },
    @{ Row = 10; Col = 5; Old = "24×56=1344"; New = "92×31=2852" },

    @{ Row = 15; Col = 1; Old = "29×83=2407"; New = "45×89=4005" },
    @{ Row = 15; Col = 2; Old = "38×96=3648"; New = "96×96=9216" },
    @{ Row = 15; Col = 3; Old = "89×27=2403"; New = "78×54=4212" },
    @{ Row = 15; Col = 4; Old = "63×59=3717"; New = "22×37=814"  },
    @{ Row = 15; Col = 5; Old = "55×19=1045"; New = "42×44=1848" },

    @{ Row = 20; Col = 1; Old = "81×98=7938"; New = "70×85=5950" },
    @{ Row = 20; Col = 2; Old = "87×77=6699"; New = "17×65=1105" },
    @{ Row = 20; Col = 3; Old = "81×19=1539"; New = "88×96=8448" },
    @{ Row = 20; Col = 4; Old = "32×93=2976"; New = "20×59=1180" },
    @{ Row = 20; Col = 5; Old = "40×26=1040"; New = "49×89=4361" }
)

$applied = 0
foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, $change.Col)
    $current = $cell.Range.Text
    if ($current.StartsWith($change.Old)) {
        $cell.Range.Text = $change.New
        $applied = $applied + 1
    } else {
        Write-Host "Skipped row" $change.Row "col" $change.Col "- unexpected content:" $current
    }
}

Write-Host "Applied" $applied "of" $changes.Count "cell updates"
